$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the new description in A9 (previously empty) - new shared string added
$ws.Range("A9").Value = "Mise à jour maquette et use case scénarios"

# Move the active selection from A9 to A10 to match the updated view state
$ws.Range("A10").Select()
